$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the "总计" (Total) roll-up sheet. The copy keeps
# accumulating the per-quarter history (gets the new 2022-Q1 row below)
# and stays named "总计"; the ORIGINAL sheet object is repurposed
# further down into the new "2022-Q1" fund-detail sheet, which is how
# the new quarter sheet ends up sitting right after "2021-Q4" and
# before "总计" in the tab order.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Copy($null, $total)
$newTotal = $wb.Worksheets.Item($total.Index + 1)
$total.Name = "2022-Q1-tmp"
$newTotal.Name = "总计"

# ---------------------------------------------------------------------
# Step 2: insert a new row 2 into the (new, copied) "总计" sheet holding
# the 2022-Q1 roll-up, pushing the existing history rows down by one.
# ---------------------------------------------------------------------
$newTotal.Rows.Item(2).Insert()
$newTotal.Rows.Item(2).ClearFormats()

# A2 needs the same bordered/bold style used by the rest of column A.
$newTotal.Range("A3").Copy()
$newTotal.Range("A2").PasteSpecial(-4122)

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 1
$newTotal.Range("D2").Value = 0.05

# Renumber the index column (A) for the rows that shifted down.
$newTotal.Range("A3").Value = 1
$newTotal.Range("A4").Value = 2
$newTotal.Range("A5").Value = 3
$newTotal.Range("A6").Value = 4
$newTotal.Range("A7").Value = 5

# ---------------------------------------------------------------------
# Step 3: turn the original "总计" sheet into the new "2022-Q1" fund
# detail sheet (same layout as the other quarterly sheets, e.g.
# "2021-Q4": fund code / name / size / position / rank columns).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1-tmp")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Apply the header/index-column style (bold, centered, thin border) to
# B1:H1 and A2, matching the other quarterly sheets - copy it from the
# "2021-Q4" sheet, which still has the original formatting intact. Do
# this BEFORE writing values so the later text-forcing step (below) is
# the last word on B2:G2's formatting.
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$styleSrc.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
$q1.Range("H2").Value = 4

# Columns B:G on row 2 look (partly) numeric but must stay stored as
# text, same as on the other quarterly sheets - force a Text number
# format before assigning so Excel doesn't silently convert them
# ("007280" -> 7280, "1.35" -> 1.35 as a Number, etc).
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "007280"
$q1.Range("C2").Value = "上投摩根日本精选股票（QDII）"
$q1.Range("D2").Value = "1.35"
$q1.Range("E2").Value = "88.71"
$q1.Range("F2").Value = "3.73"
$q1.Range("G2").Value = "0.0504"

# Restore the default (style-less) formatting on B2:G2 - only the
# stored cell *type* needs to stay text, not the visible number
# format/style index.
$q1.Range("A1").Copy()
$q1.Range("B2:G2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 4: move the new "2022-Q1" sheet so it sits right after "2021-Q4"
# and before "总计" in the tab order.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1.Move($null, $q4)
